$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metrics" (sheet1.xml) - refreshed monthly source figures
# ---------------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 286860.89
$wsMetrics.Range("B3").Value  = 246969.06
$wsMetrics.Range("B4").Value  = 90308.800000000003
$wsMetrics.Range("B5").Value  = 11360
$wsMetrics.Range("B6").Value  = 3682489.46
$wsMetrics.Range("B7").Value  = 3125683.72
$wsMetrics.Range("B8").Value  = 1053951.3599999999
$wsMetrics.Range("B9").Value  = 142048
$wsMetrics.Range("B10").Value = 32147813.259999998
$wsMetrics.Range("B11").Value = 19155553.789999999
$wsMetrics.Range("B12").Value = 11335660.25
$wsMetrics.Range("B13").Value = 1239675

# Move the active selection, as captured when the workbook was saved
$wsMetrics.Range("D11").Select()

# ---------------------------------------------------------------------------
# Sheet "today" (sheet4.xml) - refreshed daily figures; the "month so far"
# block (rows 11-22) is turned into live formulas that add the new daily
# reading onto the previous month-to-date total.
# ---------------------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("today")

$wsToday.Range("B3").Value = 16712.330000000002
$wsToday.Range("B4").Value = 14255.33
$wsToday.Range("B5").Value = 5244.38

# B6 picks up B3's number format (moves off the shared "General" style
# bucket onto the "#,##0.00" bucket used by the other daily cells).
$wsToday.Range("B6").Value = 651
$wsToday.Range("B6").NumberFormat = $wsToday.Range("B3").NumberFormat

$wsToday.Range("B11").Formula = "=270148.56+B3"
$wsToday.Range("B12").Formula = "=232713.73+B4"
$wsToday.Range("B13").Formula = "=85064.42+B5"
$wsToday.Range("B14").Formula = "=10709+B6"
$wsToday.Range("B15").Formula = '=3665777.13+$B3'
$wsToday.Range("B16").Formula = "=3111428.39+B4"
$wsToday.Range("B17").Formula = "=1048706.98+B5"
$wsToday.Range("B18").Formula = "=141397+B6"
$wsToday.Range("B19").Formula = "=32131100.93+B3"
$wsToday.Range("B20").Formula = "=19141298.46+B4"
$wsToday.Range("B21").Formula = "=11330415.87+B5"
$wsToday.Range("B22").Formula = "=1239024+B6"

# New helper cells alongside the refreshed formula rows, formatted with the
# same "#,##0.00" number format now used for the year-to-date figures
# (B18/B22 share that same format bucket).
$wsToday.Range("B18").NumberFormat = "#,##0.00_ "
$wsToday.Range("B22").NumberFormat = "#,##0.00_ "
$wsToday.Range("E11").NumberFormat = "#,##0.00_ "
$wsToday.Range("E12").NumberFormat = "#,##0.00_ "
$wsToday.Range("E13").NumberFormat = "#,##0.00_ "
$wsToday.Range("E14").NumberFormat = "#,##0.00_ "

$wb.Save()
